$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New skill rows 43-50 (column A) ---
$ws.Cells.Item(43,1).Value = "virtual machine"
$ws.Cells.Item(44,1).Value = "docker"
$ws.Cells.Item(45,1).Value = "zookeeper"
$ws.Cells.Item(46,1).Value = "go"
$ws.Cells.Item(47,1).Value = "ajax"
$ws.Cells.Item(48,1).Value = "django"
$ws.Cells.Item(49,1).Value = "perl"
$ws.Cells.Item(50,1).Value = "algorithm"

# --- New header cells E1/F1 ---
$ws.Cells.Item(1,5).Value = "testimony / certification"
$ws.Cells.Item(1,6).Value = "projects"

# --- New skill rows 51-60 (column A) ---
$ws.Cells.Item(51,1).Value = "gsm"
$ws.Cells.Item(52,1).Value = "lte"
$ws.Cells.Item(53,1).Value = "Windows GUI"
$ws.Cells.Item(54,1).Value = "Linux GUI"
$ws.Cells.Item(55,1).Value = "purify"
$ws.Cells.Item(56,1).Value = "PureCoverage"
$ws.Cells.Item(57,1).Value = "wireshark"
$ws.Cells.Item(58,1).Value = "netfilter/iptables"
$ws.Cells.Item(59,1).Value = "MAC OSX"
$ws.Cells.Item(60,1).Value = "iOS dev"

# --- Column E width (target OOXML width 23.28515625; engine quantizes to 1/6 steps, 22.5 is the closest input) ---
$ws.Columns.Item(5).ColumnWidth = 22.5

# --- Hyperlink on A56 (PureCoverage) ---
$ws.Hyperlinks.Add($ws.Range("A56"), "http://www.ltesting.net/ceshi/ceshijishu/rjcsgj/rational/purecoverage/", [Type]::Missing, [Type]::Missing, "http://www.ltesting.net/ceshi/ceshijishu/rjcsgj/rational/purecoverage/")
$ws.Range("A56").Value = "PureCoverage"
$ws.Range("A56").Style = "Normal"

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- View: scroll position + selection ---
$ws.Range("K47").Select()
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
